# Scheduled runner update: refresh market-price-derived columns (H-N)
# across the crafting-class profit sheets. Values are plain numeric
# overwrites (no formulas) pulled from the latest market snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 163265380
$ws.Range("I4").Value = 23809606
$ws.Range("K4").Value = 23809606
$ws.Range("M4").Value = -23809492

$ws.Range("H8").Value = 111111224
$ws.Range("I8").Value = 90909220
$ws.Range("J8").Value = 333333340
$ws.Range("K8").Value = 272727660
$ws.Range("L8").Value = 1000000020
$ws.Range("M8").Value = -272727521
$ws.Range("N8").Value = -1000000298

$ws.Range("H131").Value = 15294.375
$ws.Range("I131").Value = 833.3333
$ws.Range("K131").Value = 2499.9999
$ws.Range("M131").Value = 2540.0001

$ws.Range("H135").Value = 1000.0769
$ws.Range("J135").Value = 4000
$ws.Range("L135").Value = 36000
$ws.Range("N135").Value = -41070

$ws.Range("H138").Value = 3297.5647
$ws.Range("J138").Value = 3798.1667
$ws.Range("L138").Value = 11394.5001
$ws.Range("N138").Value = -21674.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 731.25
$ws.Range("I5").Value = 2066
$ws.Range("K5").Value = 2066
$ws.Range("M5").Value = -1954

$ws.Range("H61").Value = 3100.3462
$ws.Range("I61").Value = 2852.652
$ws.Range("K61").Value = 2852.652
$ws.Range("M61").Value = -2640.652

$ws.Range("H74").Value = 8441.532999999999
$ws.Range("I74").Value = 1875.2727
$ws.Range("K74").Value = 1875.2727
$ws.Range("M74").Value = -1001.2727

$ws.Range("H77").Value = 8441.532999999999
$ws.Range("I77").Value = 1875.2727
$ws.Range("K77").Value = 9376.363499999999
$ws.Range("M77").Value = -5008.363499999999

$ws.Range("H97").Value = 999.05884
$ws.Range("I97").Value = 994.3125
$ws.Range("J97").Value = 1075
$ws.Range("K97").Value = 994.3125
$ws.Range("L97").Value = 1075
$ws.Range("M97").Value = -498.3125
$ws.Range("N97").Value = -2067

$ws.Range("H132").Value = 3143.375
$ws.Range("I132").Value = 2907.3076
$ws.Range("K132").Value = 8721.9228
$ws.Range("M132").Value = -6191.9228

$ws.Range("H136").Value = 3100.3462
$ws.Range("I136").Value = 2852.652
$ws.Range("K136").Value = 8557.956
$ws.Range("M136").Value = -6007.956

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 731.25
$ws.Range("I4").Value = 2066
$ws.Range("K4").Value = 2066
$ws.Range("M4").Value = -1951

$ws.Range("H22").Value = 475.1
$ws.Range("I22").Value = 341.33334
$ws.Range("J22").Value = 675.75
$ws.Range("K22").Value = 341.33334
$ws.Range("L22").Value = 675.75
$ws.Range("M22").Value = -168.33334
$ws.Range("N22").Value = -1021.75

$ws.Range("H99").Value = 6496.5454
$ws.Range("I99").Value = 5557.875
$ws.Range("K99").Value = 5557.875
$ws.Range("M99").Value = -4059.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 50443.82
$ws.Range("I31").Value = 57221.344
$ws.Range("J31").Value = 30789
$ws.Range("K31").Value = 57221.344
$ws.Range("L31").Value = 30789
$ws.Range("M31").Value = -56926.344
$ws.Range("N31").Value = -31379

$ws.Range("H34").Value = 50443.82
$ws.Range("I34").Value = 57221.344
$ws.Range("J34").Value = 30789
$ws.Range("K34").Value = 57221.344
$ws.Range("L34").Value = 30789
$ws.Range("M34").Value = -57019.344
$ws.Range("N34").Value = -31193

$ws.Range("H132").Value = 5112.923
$ws.Range("I132").Value = 4941.222
$ws.Range("K132").Value = 14823.666
$ws.Range("M132").Value = -12293.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 9328.966
$ws.Range("I124").Value = 6112.5
$ws.Range("J124").Value = 9843.6
$ws.Range("K124").Value = 18337.5
$ws.Range("L124").Value = 29530.8
$ws.Range("M124").Value = -13427.5
$ws.Range("N124").Value = -39350.8

$ws.Range("H129").Value = 1686.125
$ws.Range("I129").Value = 762.5
$ws.Range("J129").Value = 1994
$ws.Range("K129").Value = 2287.5
$ws.Range("L129").Value = 5982
$ws.Range("M129").Value = 2712.5
$ws.Range("N129").Value = -15982

$ws.Range("H132").Value = 1322.9048
$ws.Range("I132").Value = 1199
$ws.Range("K132").Value = 10791
$ws.Range("M132").Value = -8261

$ws.Range("H137").Value = 4273.9414
$ws.Range("J137").Value = 3125
$ws.Range("L137").Value = 9375
$ws.Range("N137").Value = -19575

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8836.058999999999
$ws.Range("I70").Value = 8619.272000000001
$ws.Range("K70").Value = 8619.272000000001
$ws.Range("M70").Value = -8349.272000000001

$ws.Range("H73").Value = 8836.058999999999
$ws.Range("I73").Value = 8619.272000000001
$ws.Range("K73").Value = 8619.272000000001
$ws.Range("M73").Value = -7683.272000000001

$ws.Range("H126").Value = 14306.308
$ws.Range("I126").Value = 18472.422
$ws.Range("K126").Value = 55417.266
$ws.Range("M126").Value = -52947.266

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4040.7778
$ws.Range("I22").Value = 1406.125
$ws.Range("J22").Value = 6148.5
$ws.Range("K22").Value = 1406.125
$ws.Range("L22").Value = 6148.5
$ws.Range("M22").Value = -1111.125
$ws.Range("N22").Value = -6738.5

$ws.Range("H27").Value = 4040.7778
$ws.Range("I27").Value = 1406.125
$ws.Range("J27").Value = 6148.5
$ws.Range("K27").Value = 1406.125
$ws.Range("L27").Value = 6148.5
$ws.Range("M27").Value = -1299.125
$ws.Range("N27").Value = -6362.5

$ws.Range("H55").Value = 189.42857
$ws.Range("I55").Value = 212
$ws.Range("K55").Value = 212
$ws.Range("M55").Value = -39

$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0

$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0

$ws.Range("H136").Value = 5999.316
$ws.Range("I136").Value = 5599.2
$ws.Range("J136").Value = 7499.75
$ws.Range("K136").Value = 16797.6
$ws.Range("L136").Value = 22499.25
$ws.Range("M136").Value = -14247.6
$ws.Range("N136").Value = -27599.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1301.9546
$ws.Range("I100").Value = 1086.3334
$ws.Range("K100").Value = 2172.6668
$ws.Range("M100").Value = -1631.6668

$ws.Range("H107").Value = 1762.3889
$ws.Range("I107").Value = 1789.2
$ws.Range("J107").Value = 1752.0769
$ws.Range("K107").Value = 5367.6
$ws.Range("L107").Value = 5256.2307
$ws.Range("M107").Value = -3447.6
$ws.Range("N107").Value = -9096.2307

$ws.Range("H126").Value = 3645.625
$ws.Range("I126").Value = 3352.1177
$ws.Range("K126").Value = 10056.3531
$ws.Range("M126").Value = -7586.3531

# Rows 88 and 91 on LTW no longer carry a LeveProfitNQ figure at all
# (the source sheet drops the cell rather than zeroing it out).
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M88").ClearContents()
$ws.Range("M91").ClearContents()
